# Update "Advies" sheet: new timestamp + refreshed advice table (rows 3-17)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Laatst bijgewerkt" timestamp (row 2, col A)
$ws.Cells.Item(2, 1).Value = "Laatst bijgewerkt: 2025-09-05 18:53:17"

# 2) Remove all existing hyperlinks up-front so re-adding them below never
#    stacks duplicate hyperlink objects on top of the old ones.
$ws.Hyperlinks.Delete()

$rows = @(
  @{ Row=3; A="Chicago Fire vs New England Revolution"; B="totaal aantal schoten op doel"; C="wedstrijd"; D="meer dan 10.5"; E="onecasino"; F=2.73; G="minder dan 10.5"; H="starcasino"; I=1.75; J="1=59, 2=91"; K="€9.25"; L=6.23; M="https://sports.onecasino.nl/#/event/10025497"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=11998073" },
  @{ Row=4; A="Duitsland vs Noord-Ierland"; B="totaal aantal schoten op doel"; C="duitsland"; D="meer dan 7.5"; E="bingoal"; F=2.1; G="minder dan 7.5"; H="jacks"; I=2.08; J="1=75, 2=75"; K="€6.0"; L=4.3; M="https://www.bingoal.nl/sports/#event=1023224852&betoffer=2552162167&outcome=3865656359"; N="https://jacks.nl/sports/event/1023224852#event/1023224852" },
  @{ Row=5; A="Duitsland vs Noord-Ierland"; B="totaal aantal schoten op doel"; C="duitsland"; D="meer dan 8.5"; E="bingoal"; F=2.85; G="minder dan 8.5"; H="jacks"; I=1.64; J="1=55, 2=95"; K="€5.8"; L=3.94; M="https://www.bingoal.nl/sports/#event=1023224852&betoffer=2552162170&outcome=3865656366"; N="https://jacks.nl/sports/event/1023224852#event/1023224852" },
  @{ Row=6; A="Montenegro vs Tsjechië"; B="totaal aantal schoten op doel"; C="montenegro"; D="meer dan 3.5"; E="bingoal"; F=2.6; G="minder dan 3.5"; H="starcasino"; I=1.7273; J="1=60, 2=90"; K="€5.46"; L=3.64; M="https://www.bingoal.nl/sports/#event=1023168639&betoffer=2549626618&outcome=3856355987"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=13549442" },
  @{ Row=7; A="Duitsland vs Noord-Ierland"; B="totaal aantal schoten op doel"; C="duitsland"; D="meer dan 6.5"; E="bingoal"; F=1.61; G="minder dan 6.5"; H="jacks"; I=2.88; J="1=96, 2=54"; K="€4.56"; L=3.17; M="https://www.bingoal.nl/sports/#event=1023224852&betoffer=2552162169&outcome=3865656364"; N="https://jacks.nl/sports/event/1023224852#event/1023224852" },
  @{ Row=8; A="Montenegro vs Tsjechië"; B="totaal aantal schoten op doel"; C="montenegro"; D="meer dan 3.5"; E="jacks"; F=2.55; G="minder dan 3.5"; H="starcasino"; I=1.7273; J="1=61, 2=89"; K="€3.73"; L=2.89; M="https://jacks.nl/sports/event/1023168639#event/1023168639"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=13549442" },
  @{ Row=9; A="Oekraïne vs Frankrijk"; B="totaal aantal schoten"; C="kylian mbappé"; D="meer dan 4.5"; E="toto"; F=2.09; G="minder dan 4.5"; H="jacks"; I=2.02; J="1=74, 2=76"; K="€3.52"; L=2.65; M="https://sport.toto.nl/wedden/wedstrijd/8590797"; N="https://jacks.nl/sports/event/1023224887#event/1023224887" },
  @{ Row=10; A="Oekraïne vs Frankrijk"; B="totaal aantal schoten op doel"; C="oekraïne"; D="meer dan 3.5"; E="toto"; F=3.5; G="minder dan 3.5"; H="starcasino"; I=1.4445; J="1=44, 2=106"; K="€3.12"; L=2.2; M="https://sport.toto.nl/wedden/wedstrijd/8590797"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12623890" },
  @{ Row=11; A="IJsland vs Azerbeidzjan"; B="totaal aantal schoten op doel"; C="azerbeidzjan"; D="meer dan 3.5"; E="bingoal"; F=2.33; G="minder dan 3.5"; H="starcasino"; I=1.8182; J="1=66, 2=84"; K="€2.73"; L=2.08; M="https://www.bingoal.nl/sports/#event=1022335870&betoffer=2549614664&outcome=3856406727"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12623891" },
  @{ Row=12; A="Montenegro vs Tsjechië"; B="totaal aantal schoten op doel"; C="montenegro"; D="meer dan 2.5"; E="bingoal"; F=1.65; G="minder dan 2.5"; H="starcasino"; I=2.6667; J="1=93, 2=57"; K="€2.0"; L=1.89; M="https://www.bingoal.nl/sports/#event=1023168639&betoffer=2549626614&outcome=3856355980"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=13549442" },
  @{ Row=13; A="Montenegro vs Tsjechië"; B="totaal aantal schoten op doel"; C="montenegro"; D="meer dan 2.5"; E="jacks"; F=1.64; G="minder dan 2.5"; H="starcasino"; I=2.6667; J="1=93, 2=57"; K="€2.0"; L=1.52; M="https://jacks.nl/sports/event/1023168639#event/1023168639"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=13549442" },
  @{ Row=14; A="Oekraïne vs Frankrijk"; B="totaal aantal schoten op doel"; C="wedstrijd"; D="meer dan 8.5"; E="toto"; F=2.05; G="minder dan 8.5"; H="starcasino"; I=2; J="1=74, 2=76"; K="€1.7"; L=1.22; M="https://sport.toto.nl/wedden/wedstrijd/8590797"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12623890" },
  @{ Row=15; A="Oekraïne vs Frankrijk"; B="totaal aantal schoten op doel"; C="oekraïne"; D="meer dan 2.5"; E="toto"; F=2.05; G="minder dan 2.5"; H="starcasino"; I=2; J="1=74, 2=76"; K="€1.7"; L=1.22; M="https://sport.toto.nl/wedden/wedstrijd/8590797"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12623890" },
  @{ Row=16; A="Duitsland vs Noord-Ierland"; B="totaal aantal schoten op doel"; C="wedstrijd"; D="meer dan 10.5"; E="toto"; F=2.5; G="minder dan 10.5"; H="onecasino"; I=1.7; J="1=61, 2=89"; K="€1.3"; L=1.18; M="https://sport.toto.nl/wedden/wedstrijd/8778584"; N="https://sports.onecasino.nl/#/event/10028349" },
  @{ Row=17; A="Oekraïne vs Frankrijk"; B="totaal aantal schoten op doel"; C="wedstrijd"; D="meer dan 9.5"; E="toto"; F=2.75; G="minder dan 9.5"; H="starcasino"; I=1.6; J="1=55, 2=95"; K="€1.25"; L=1.14; M="https://sport.toto.nl/wedden/wedstrijd/8590797"; N="https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12623890" }
)

# 3) Write the full advice table (header row 1 is unchanged)
foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J

    # Column K ("Beste inzet") holds euro-amount TEXT like "€6.0" / "€2.0" -
    # force text format first so Excel doesn't silently coerce it to a
    # currency number (which would also eat the trailing zero).
    $cK = $ws.Cells.Item($row, 11)
    $cK.NumberFormat = "@"
    $cK.Value = $r.K

    $ws.Cells.Item($row, 12).Value = $r.L

    # Columns M/N ("Link 1"/"Link 2") are hyperlinks. The cell text is the
    # full URL (base + "#" + fragment); Address is the part before "#" and
    # SubAddress (location) is the part after it, matching how Excel stores
    # hyperlink Target/location in the OOXML.
    $mText = $r.M
    $mIdx = $mText.IndexOf("#")
    if ($mIdx -eq -1) {
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 13), $mText, "", "", $mText) | Out-Null
    } else {
        $mBase = $mText.Substring(0, $mIdx)
        $mFrag = $mText.Substring($mIdx + 1)
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 13), $mBase, $mFrag, "", $mText) | Out-Null
    }

    $nText = $r.N
    $nIdx = $nText.IndexOf("#")
    if ($nIdx -eq -1) {
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 14), $nText, "", "", $nText) | Out-Null
    } else {
        $nBase = $nText.Substring(0, $nIdx)
        $nFrag = $nText.Substring($nIdx + 1)
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 14), $nBase, $nFrag, "", $nText) | Out-Null
    }
}
